# edit.ps1 - apply the text-merge fixes described by the diff.
#
# Summary of changes:
#  1. Slide 11 ("Budget and costs"), TextBox 2: "Miscellanous" + " " -> "Miscellaneous "
#     (fix typo and merge the two runs into a single run).
#  2. Slide 13 ("Sprint planning"), TextBox 2: "Research " + "Symposium: 4/7/2016"
#     -> "Research Symposium: 4/7/2016" (merge the two runs into a single run).
#  3. Slide 4 ("The Cluster"), TextBox 2: "Snow White and the seven " + "dwarfs"
#     -> "Snow White and the seven dwarfs" (merge the two runs into a single run).
#
# All merges are done surgically via TextRange.Characters() ranges so the
# surrounding paragraph/run structure (and their formatting) is otherwise left
# untouched: the run that already carries the "clean" rPr (no spurious
# err="1" / correct formatting) is expanded to hold the full merged text and
# the other (now redundant) run is cleared out so it disappears from the XML.

$p = $ppt.ActivePresentation

# --- 1) Slide 11: "Miscellanous" + " " -> "Miscellaneous " ------------------
$s11 = $p.Slides.Item(11)
$shape11 = $s11.Shapes.Item(2)
$tr11 = $shape11.TextFrame.TextRange
$para11 = $tr11.Paragraphs(25, 1)
# run2 (" ") keeps clean formatting (no err="1") -> grow it to the full fixed text
$run2_11 = $para11.Characters(13, 1)
$run2_11.Text = "Miscellaneous "
# run1 ("Miscellanous") becomes redundant -> clear it out
$run1_11 = $para11.Characters(1, 12)
$run1_11.Text = ""

# --- 2) Slide 13: "Research " + "Symposium: 4/7/2016" ----------------------
$s13 = $p.Slides.Item(13)
$shape13 = $s13.Shapes.Item(2)
$tr13 = $shape13.TextFrame.TextRange
$para13 = $tr13.Paragraphs(15, 1)
$full13 = $para13.Characters(1, $para13.Text.Length)
$full13.Text = "Research Symposium: 4/7/2016"

# --- 3) Slide 4: "Snow White and the seven " + "dwarfs" --------------------
$s4 = $p.Slides.Item(4)
$shape4 = $s4.Shapes.Item(2)
$tr4 = $shape4.TextFrame.TextRange
$para4 = $tr4.Paragraphs(1, 1)
$full4 = $para4.Characters(1, $para4.Text.Length)
$full4.Text = "Snow White and the seven dwarfs"
